# Apply updated loading_percent results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.36711817794491
$ws.Cells.Item(2, 4).Value = 6.198157109497284
$ws.Cells.Item(2, 5).Value = 18.66654321818146
$ws.Cells.Item(2, 6).Value = 32.21748607642428
$ws.Cells.Item(2, 7).Value = 42.7734385296047
$ws.Cells.Item(2, 8).Value = 17.52470279179087
$ws.Cells.Item(2, 11).Value = 9.418769110415182
$ws.Cells.Item(2, 12).Value = 8.548674888827144
$ws.Cells.Item(2, 13).Value = 15.26759374646476
$ws.Cells.Item(3, 2).Value = 16.30682127409546
$ws.Cells.Item(3, 4).Value = 6.196388403747609
$ws.Cells.Item(3, 5).Value = 18.68875650849771
$ws.Cells.Item(3, 6).Value = 32.08962860626016
$ws.Cells.Item(3, 7).Value = 42.51446890057819
$ws.Cells.Item(3, 8).Value = 17.53131773892376
$ws.Cells.Item(3, 11).Value = 9.081209032844825
$ws.Cells.Item(3, 12).Value = 8.5423535061928
$ws.Cells.Item(3, 13).Value = 15.26455463064141
$ws.Cells.Item(4, 2).Value = 16.27314975577361
$ws.Cells.Item(4, 4).Value = 6.195186779575255
$ws.Cells.Item(4, 5).Value = 18.70323983631898
$ws.Cells.Item(4, 6).Value = 32.01848840389974
$ws.Cells.Item(4, 7).Value = 42.36576599999086
$ws.Cells.Item(4, 8).Value = 17.53834962116805
$ws.Cells.Item(4, 11).Value = 8.865074799309191
$ws.Cells.Item(4, 12).Value = 8.539960626491414
$ws.Cells.Item(4, 13).Value = 15.26519370735551
$ws.Cells.Item(5, 2).Value = 16.26028189551645
$ws.Cells.Item(5, 4).Value = 6.194668129143385
$ws.Cells.Item(5, 5).Value = 18.70935458614172
$ws.Cells.Item(5, 6).Value = 31.99136855293487
$ws.Cells.Item(5, 7).Value = 42.30780566756822
$ws.Cells.Item(5, 8).Value = 17.54196171571396
$ws.Cells.Item(5, 11).Value = 8.774834099524011
$ws.Cells.Item(5, 12).Value = 8.539361216010384
$ws.Cells.Item(5, 13).Value = 15.26608585141068
$ws.Cells.Item(6, 2).Value = 16.25819704949371
$ws.Cells.Item(6, 4).Value = 6.19458026029442
$ws.Cells.Item(6, 5).Value = 18.71038279189995
$ws.Cells.Item(6, 6).Value = 31.98697877696377
$ws.Cells.Item(6, 7).Value = 42.29834193572876
$ws.Cells.Item(6, 8).Value = 17.54260657527773
$ws.Cells.Item(6, 11).Value = 8.759721086952412
$ws.Cells.Item(6, 12).Value = 8.539284413745115
$ws.Cells.Item(6, 13).Value = 15.2662721924279
$ws.Cells.Item(7, 2).Value = 16.27297274560896
$ws.Cells.Item(7, 4).Value = 6.195179902090106
$ws.Cells.Item(7, 5).Value = 18.70332144033926
$ws.Cells.Item(7, 6).Value = 32.01811506067621
$ws.Cells.Item(7, 7).Value = 42.36497359121567
$ws.Cells.Item(7, 8).Value = 17.5383953131089
$ws.Cells.Item(7, 11).Value = 8.863866450845286
$ws.Cells.Item(7, 12).Value = 8.539951019681373
$ws.Cells.Item(7, 13).Value = 15.26520317916859
$ws.Cells.Item(8, 2).Value = 16.34563914049912
$ws.Cells.Item(8, 4).Value = 6.197571316868835
$ws.Cells.Item(8, 5).Value = 18.67402735727475
$ws.Cells.Item(8, 6).Value = 32.1718840452597
$ws.Cells.Item(8, 7).Value = 42.68203582408191
$ws.Cells.Item(8, 8).Value = 17.52636700247197
$ws.Cells.Item(8, 11).Value = 9.304256936902972
$ws.Cells.Item(8, 12).Value = 8.546187250318454
$ws.Cells.Item(8, 13).Value = 15.26602706840377
$ws.Cells.Item(9, 2).Value = 16.5142036234835
$ws.Cells.Item(9, 4).Value = 6.201339155432288
$ws.Cells.Item(9, 5).Value = 18.62326359975063
$ws.Cells.Item(9, 6).Value = 32.53094091932436
$ws.Cells.Item(9, 7).Value = 43.38333450060304
$ws.Cells.Item(9, 8).Value = 17.52635476053406
$ws.Cells.Item(9, 11).Value = 10.09495555091429
$ws.Cells.Item(9, 12).Value = 8.570157104897703
$ws.Cells.Item(9, 13).Value = 15.28743207138861
$ws.Cells.Item(10, 2).Value = 16.65317759625137
$ws.Cells.Item(10, 4).Value = 6.203539869742071
$ws.Cells.Item(10, 5).Value = 18.5900169918466
$ws.Cells.Item(10, 6).Value = 32.82841697122566
$ws.Cells.Item(10, 7).Value = 43.9438003835952
$ws.Cells.Item(10, 8).Value = 17.54071554910849
$ws.Cells.Item(10, 11).Value = 10.62857434829057
$ws.Cells.Item(10, 12).Value = 8.59482074583017
$ws.Cells.Item(10, 13).Value = 15.31509163486459
$ws.Cells.Item(11, 2).Value = 16.71950569361325
$ws.Cells.Item(11, 4).Value = 6.204416656695607
$ws.Cells.Item(11, 5).Value = 18.5757665250345
$ws.Cells.Item(11, 6).Value = 32.97071739685269
$ws.Cells.Item(11, 7).Value = 44.20781706930433
$ws.Cells.Item(11, 8).Value = 17.55036281909575
$ws.Cells.Item(11, 11).Value = 10.8605827542563
$ws.Cells.Item(11, 12).Value = 8.607546023513441
$ws.Cells.Item(11, 13).Value = 15.33023259263548
$ws.Cells.Item(12, 2).Value = 16.74505296131385
$ws.Cells.Item(12, 4).Value = 6.204730704704466
$ws.Cells.Item(12, 5).Value = 18.57049549340372
$ws.Cells.Item(12, 6).Value = 33.02557412059851
$ws.Cells.Item(12, 7).Value = 44.30902602549376
$ws.Cells.Item(12, 8).Value = 17.55446276321283
$ws.Cells.Item(12, 11).Value = 10.94686165204207
$ws.Cells.Item(12, 12).Value = 8.612578687818848
$ws.Cells.Item(12, 13).Value = 15.33633067655626
$ws.Cells.Item(13, 2).Value = 16.739532037431
$ws.Cells.Item(13, 4).Value = 6.204663869623152
$ws.Cells.Item(13, 5).Value = 18.57162513602002
$ws.Cells.Item(13, 6).Value = 33.01371711359614
$ws.Cells.Item(13, 7).Value = 44.28717524839378
$ws.Cells.Item(13, 8).Value = 17.55355991945457
$ws.Cells.Item(13, 11).Value = 10.92835059756298
$ws.Cells.Item(13, 12).Value = 8.611485345299908
$ws.Cells.Item(13, 13).Value = 15.33500118907051
$ws.Cells.Item(14, 2).Value = 16.72159895777964
$ws.Cells.Item(14, 4).Value = 6.204442854331297
$ws.Cells.Item(14, 5).Value = 18.57533036442611
$ws.Cells.Item(14, 6).Value = 32.97521123304053
$ws.Cells.Item(14, 7).Value = 44.21611932410046
$ws.Cells.Item(14, 8).Value = 17.55069117830039
$ws.Cells.Item(14, 11).Value = 10.86771276474095
$ws.Cells.Item(14, 12).Value = 8.607955793200121
$ws.Cells.Item(14, 13).Value = 15.33072700021212
$ws.Cells.Item(15, 2).Value = 16.71066995358421
$ws.Cells.Item(15, 4).Value = 6.204305132838455
$ws.Cells.Item(15, 5).Value = 18.57761623389285
$ws.Cells.Item(15, 6).Value = 32.95175067470318
$ws.Cells.Item(15, 7).Value = 44.17275375540813
$ws.Cells.Item(15, 8).Value = 17.54899213071531
$ws.Cells.Item(15, 11).Value = 10.83036397071081
$ws.Cells.Item(15, 12).Value = 8.605821612249716
$ws.Cells.Item(15, 13).Value = 15.32815630409088
$ws.Cells.Item(16, 2).Value = 16.64890403227313
$ws.Cells.Item(16, 4).Value = 6.203480059746816
$ws.Cells.Item(16, 5).Value = 18.59096584553892
$ws.Cells.Item(16, 6).Value = 32.81925495186165
$ws.Cells.Item(16, 7).Value = 43.92672253616061
$ws.Cells.Item(16, 8).Value = 17.54014765918387
$ws.Cells.Item(16, 11).Value = 10.61319300096701
$ws.Cells.Item(16, 12).Value = 8.594019196127629
$ws.Cells.Item(16, 13).Value = 15.31415333865286
$ws.Cells.Item(17, 2).Value = 16.61179706183517
$ws.Cells.Item(17, 4).Value = 6.202941992810689
$ws.Cells.Item(17, 5).Value = 18.59937892373742
$ws.Cells.Item(17, 6).Value = 32.73973748321629
$ws.Cells.Item(17, 7).Value = 43.77806067958595
$ws.Cells.Item(17, 8).Value = 17.53551884884286
$ws.Cells.Item(17, 11).Value = 10.47719000891604
$ws.Cells.Item(17, 12).Value = 8.587162621467206
$ws.Cells.Item(17, 13).Value = 15.30621594782576
$ws.Cells.Item(18, 2).Value = 16.59074771279775
$ws.Cells.Item(18, 4).Value = 6.202620813879099
$ws.Cells.Item(18, 5).Value = 18.60430014770555
$ws.Cells.Item(18, 6).Value = 32.69465993286126
$ws.Cells.Item(18, 7).Value = 43.6934113171376
$ws.Cells.Item(18, 8).Value = 17.53314971494131
$ws.Cells.Item(18, 11).Value = 10.39795569842495
$ws.Cells.Item(18, 12).Value = 8.583360791695403
$ws.Cells.Item(18, 13).Value = 15.30189160656583
$ws.Cells.Item(19, 2).Value = 16.58367168766482
$ws.Cells.Item(19, 4).Value = 6.202510062150519
$ws.Cells.Item(19, 5).Value = 18.60598052453267
$ws.Cells.Item(19, 6).Value = 32.6795115615404
$ws.Cells.Item(19, 7).Value = 43.66489980735186
$ws.Cells.Item(19, 8).Value = 17.5323979597866
$ws.Cells.Item(19, 11).Value = 10.37095613410689
$ws.Cells.Item(19, 12).Value = 8.582098004327424
$ws.Cells.Item(19, 13).Value = 15.30046895543105
$ws.Cells.Item(20, 2).Value = 16.6157168944553
$ws.Cells.Item(20, 4).Value = 6.203000482232817
$ws.Cells.Item(20, 5).Value = 18.59847482753297
$ws.Cells.Item(20, 6).Value = 32.74813430128678
$ws.Cells.Item(20, 7).Value = 43.79379779673054
$ws.Cells.Item(20, 8).Value = 17.53598125200475
$ws.Cells.Item(20, 11).Value = 10.4917724451227
$ws.Cells.Item(20, 12).Value = 8.587877848850828
$ws.Cells.Item(20, 13).Value = 15.30703597388113
$ws.Cells.Item(21, 2).Value = 16.72685479738232
$ws.Cells.Item(21, 4).Value = 6.204508260372044
$ws.Cells.Item(21, 5).Value = 18.5742386505685
$ws.Cells.Item(21, 6).Value = 32.98649527002028
$ws.Cells.Item(21, 7).Value = 44.2369573252343
$ws.Cells.Item(21, 8).Value = 17.5515216833852
$ws.Cells.Item(21, 11).Value = 10.88556662676831
$ws.Cells.Item(21, 12).Value = 8.608986725606838
$ws.Cells.Item(21, 13).Value = 15.33197256805105
$ws.Cells.Item(22, 2).Value = 16.80198949138088
$ws.Cells.Item(22, 4).Value = 6.205388883452594
$ws.Cells.Item(22, 5).Value = 18.55912917536648
$ws.Cells.Item(22, 6).Value = 33.14791916817796
$ws.Cells.Item(22, 7).Value = 44.53373159127383
$ws.Cells.Item(22, 8).Value = 17.56428141739293
$ws.Cells.Item(22, 11).Value = 11.13372588693541
$ws.Cells.Item(22, 12).Value = 8.624027879508368
$ws.Cells.Item(22, 13).Value = 15.35039308661352
$ws.Cells.Item(23, 2).Value = 16.76166566253198
$ws.Cells.Item(23, 4).Value = 6.204928497543732
$ws.Cells.Item(23, 5).Value = 18.5671266680415
$ws.Cells.Item(23, 6).Value = 33.06125927942717
$ws.Cells.Item(23, 7).Value = 44.37470799289276
$ws.Cells.Item(23, 8).Value = 17.55723356706979
$ws.Cells.Item(23, 11).Value = 11.0021310175262
$ws.Cells.Item(23, 12).Value = 8.615887110525957
$ws.Cells.Item(23, 13).Value = 15.34036863878948
$ws.Cells.Item(24, 2).Value = 16.61394385240928
$ws.Cells.Item(24, 4).Value = 6.202974076005622
$ws.Cells.Item(24, 5).Value = 18.59888330642854
$ws.Cells.Item(24, 6).Value = 32.74433610964102
$ws.Cells.Item(24, 7).Value = 43.7866804915874
$ws.Cells.Item(24, 8).Value = 17.53577128973155
$ws.Cells.Item(24, 11).Value = 10.48518297621321
$ws.Cells.Item(24, 12).Value = 8.587554057933392
$ws.Cells.Item(24, 13).Value = 15.30666449547171
$ws.Cells.Item(25, 2).Value = 16.46588693145116
$ws.Cells.Item(25, 4).Value = 6.200419123285435
$ws.Cells.Item(25, 5).Value = 18.63628369135143
$ws.Cells.Item(25, 6).Value = 32.42778304912393
$ws.Cells.Item(25, 7).Value = 43.18542128348601
$ws.Cells.Item(25, 8).Value = 17.52383329449234
$ws.Cells.Item(25, 11).Value = 9.889177009973226
$ws.Cells.Item(25, 12).Value = 8.562425678329426
$ws.Cells.Item(25, 13).Value = 15.27953528613736
